$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1226
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 228433
$ws.Range("I6").Value = 1333506.6
$ws.Range("J6").Value = 53947.684
$ws.Range("K6").Value = 4000519.8
$ws.Range("L6").Value = 161843.052
$ws.Range("M6").Value = -4000407.8
$ws.Range("N6").Value = -162067.052
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 466.31818
$ws.Range("I41").Value = 289.08334
$ws.Range("J41").Value = 679
$ws.Range("K41").Value = 289.08334
$ws.Range("L41").Value = 679
$ws.Range("M41").Value = 150.91666
$ws.Range("N41").Value = -1559
# Row 45 (Leve Item ID 4585)
$ws.Range("H45").Value = 250
$ws.Range("J45").Value = 250
$ws.Range("L45").Value = 750
$ws.Range("N45").Value = -1134
# Row 53 (Leve Item ID 5479)
$ws.Range("H53").Value = 2434.0557
$ws.Range("I53").Value = 2983.7144
$ws.Range("K53").Value = 2983.7144
$ws.Range("M53").Value = -2346.7144
# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 2683.3333
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15300
# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 1265.1765
$ws.Range("I92").Value = 1218.9375
$ws.Range("J92").Value = 2005
$ws.Range("K92").Value = 1218.9375
$ws.Range("L92").Value = 2005
$ws.Range("M92").Value = 29.0625
$ws.Range("N92").Value = -4501
# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 36720.688
$ws.Range("I98").Value = 38635.332
$ws.Range("K98").Value = 38635.332
$ws.Range("M98").Value = -37137.332
# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 6502705
$ws.Range("I106").Value = 7266405.5
$ws.Range("K106").Value = 7266405.5
$ws.Range("M106").Value = -7265774.5
# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 2918.25
$ws.Range("J112").Value = 1989.4
$ws.Range("L112").Value = 5968.200000000001
$ws.Range("N112").Value = -8184.200000000001
# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 36720.688
$ws.Range("I122").Value = 38635.332
$ws.Range("K122").Value = 115905.996
$ws.Range("M122").Value = -113455.996
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1420319.4
$ws.Range("I137").Value = 2321113.2
$ws.Range("J137").Value = 4785.857
$ws.Range("K137").Value = 6963339.600000001
$ws.Range("L137").Value = 14357.571
$ws.Range("M137").Value = -6960789.600000001
$ws.Range("N137").Value = -19457.571
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 4975.7593
$ws.Range("I138").Value = 2341.8572
$ws.Range("J138").Value = 5543.0615
$ws.Range("K138").Value = 7025.571599999999
$ws.Range("L138").Value = 16629.1845
$ws.Range("M138").Value = -1885.571599999999
$ws.Range("N138").Value = -26909.1845

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 25 (Leve Item ID 2471)
$ws.Range("H25").Value = 6398.2
$ws.Range("I25").Value = 3997.3333
$ws.Range("J25").Value = 9999.5
$ws.Range("K25").Value = 3997.3333
$ws.Range("L25").Value = 9999.5
$ws.Range("M25").Value = -3595.3333
$ws.Range("N25").Value = -10803.5
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 6548.4736
$ws.Range("I32").Value = 6212.5093
$ws.Range("K32").Value = 6212.5093
$ws.Range("M32").Value = -5925.5093
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 2090.9395
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 2090.9395
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 942176
$ws.Range("J122").Value = 3338899.8
$ws.Range("L122").Value = 10016699.4
$ws.Range("N122").Value = -10021599.4

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 37 (Leve Item ID 2485)
$ws.Range("H37").Value = 20000
$ws.Range("I37").Value = 20000
$ws.Range("K37").Value = 20000
$ws.Range("M37").Value = -19863
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 9043.777
$ws.Range("I105").Value = 9445.23
$ws.Range("K105").Value = 9445.23
$ws.Range("M105").Value = -7698.23
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1826.6875
$ws.Range("I134").Value = 1237.7858
$ws.Range("J134").Value = 5949
$ws.Range("K134").Value = 3713.3574
$ws.Range("L134").Value = 17847
$ws.Range("M134").Value = -1178.3574
$ws.Range("N134").Value = -22917

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 277.3684
$ws.Range("I7").Value = 226.46666
$ws.Range("J7").Value = 468.25
$ws.Range("K7").Value = 226.46666
$ws.Range("L7").Value = 468.25
$ws.Range("M7").Value = -113.46666
$ws.Range("N7").Value = -694.25
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 3892.7144
$ws.Range("I16").Value = 5750
$ws.Range("K16").Value = 5750
$ws.Range("M16").Value = -5463
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2900.5
$ws.Range("I31").Value = 1870.75
$ws.Range("J31").Value = 6333
$ws.Range("K31").Value = 1870.75
$ws.Range("L31").Value = 6333
$ws.Range("M31").Value = -1575.75
$ws.Range("N31").Value = -6923
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2900.5
$ws.Range("I34").Value = 1870.75
$ws.Range("J34").Value = 6333
$ws.Range("K34").Value = 1870.75
$ws.Range("L34").Value = 6333
$ws.Range("M34").Value = -1668.75
$ws.Range("N34").Value = -6737
# Row 103 (Leve Item ID 19558)
$ws.Range("H103").Value = 29929.072
$ws.Range("I103").Value = 18135.363
$ws.Range("J103").Value = 73172.664
$ws.Range("K103").Value = 18135.363
$ws.Range("L103").Value = 73172.664
$ws.Range("M103").Value = -16963.363
$ws.Range("N103").Value = -75516.664
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 3892.7144
$ws.Range("I113").Value = 5750
$ws.Range("K113").Value = 5750
$ws.Range("M113").Value = -3580
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 6820.778
$ws.Range("I132").Value = 8821.462
$ws.Range("K132").Value = 26464.386
$ws.Range("M132").Value = -23934.386

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 8 (Leve Item ID 16734)
$ws.Range("H8").Value = 455.83334
$ws.Range("I8").Value = 455.83334
$ws.Range("K8").Value = 1367.50002
$ws.Range("M8").Value = -1228.50002
# Row 56 (Leve Item ID 10146)
$ws.Range("H56").Value = 7942.6553
$ws.Range("I56").Value = 7942.6553
$ws.Range("K56").Value = 7942.6553
$ws.Range("M56").Value = -7412.6553
# Row 87 (Leve Item ID 12864)
$ws.Range("H87").Value = 17120.54
$ws.Range("I87").Value = 12927.25
$ws.Range("J87").Value = 18984.223
$ws.Range("K87").Value = 38781.75
$ws.Range("L87").Value = 56952.66900000001
$ws.Range("M87").Value = -37533.75
$ws.Range("N87").Value = -59448.66900000001
# Row 90 (Leve Item ID 12864)
$ws.Range("H90").Value = 17120.54
$ws.Range("I90").Value = 12927.25
$ws.Range("J90").Value = 18984.223
$ws.Range("K90").Value = 116345.25
$ws.Range("L90").Value = 170858.007
$ws.Range("M90").Value = -110105.25
$ws.Range("N90").Value = -183338.007
# Row 117 (Leve Item ID 27870)
$ws.Range("H117").Value = 1015.5
$ws.Range("J117").Value = 2000
$ws.Range("L117").Value = 6000
$ws.Range("N117").Value = -12884

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 8813.218999999999
$ws.Range("I70").Value = 7229.875
$ws.Range("J70").Value = 10396.5625
$ws.Range("K70").Value = 7229.875
$ws.Range("L70").Value = 10396.5625
$ws.Range("M70").Value = -6959.875
$ws.Range("N70").Value = -10936.5625
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 8813.218999999999
$ws.Range("I73").Value = 7229.875
$ws.Range("J73").Value = 10396.5625
$ws.Range("K73").Value = 7229.875
$ws.Range("L73").Value = 10396.5625
$ws.Range("M73").Value = -6293.875
$ws.Range("N73").Value = -12268.5625
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 13607.467
$ws.Range("I80").Value = 16313.7
$ws.Range("K80").Value = 16313.7
$ws.Range("M80").Value = -15315.7
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 13607.467
$ws.Range("I83").Value = 16313.7
$ws.Range("K83").Value = 81568.5
$ws.Range("M83").Value = -76576.5
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 8158.1177
$ws.Range("I97").Value = 8045.933
$ws.Range("J97").Value = 8999.5
$ws.Range("K97").Value = 8045.933
$ws.Range("L97").Value = 8999.5
$ws.Range("M97").Value = -7549.933
$ws.Range("N97").Value = -9991.5
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 15014.889
$ws.Range("I113").Value = 20022.334
$ws.Range("K113").Value = 20022.334
$ws.Range("M113").Value = -17852.334
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 16863.48
$ws.Range("I122").Value = 11439.333
$ws.Range("J122").Value = 24999.7
$ws.Range("K122").Value = 34317.999
$ws.Range("L122").Value = 74999.10000000001
$ws.Range("M122").Value = -31867.999
$ws.Range("N122").Value = -79899.10000000001
# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 43500
$ws.Range("J123").Value = 43500
$ws.Range("L123").Value = 43500
$ws.Range("N123").Value = -48400

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 6711.875
$ws.Range("I16").Value = 7308.05
$ws.Range("K16").Value = 7308.05
$ws.Range("M16").Value = -7138.05
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 910.4167
$ws.Range("J22").Value = 970.3333
$ws.Range("L22").Value = 970.3333
$ws.Range("N22").Value = -1560.3333
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 910.4167
$ws.Range("J27").Value = 970.3333
$ws.Range("L27").Value = 970.3333
$ws.Range("N27").Value = -1184.3333
# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 25204.955
$ws.Range("I61").Value = 1973.9333
$ws.Range("J61").Value = 74985.71000000001
$ws.Range("K61").Value = 1973.9333
$ws.Range("L61").Value = 74985.71000000001
$ws.Range("M61").Value = -1771.9333
$ws.Range("N61").Value = -75389.71000000001
# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 25204.955
$ws.Range("I113").Value = 1973.9333
$ws.Range("J113").Value = 74985.71000000001
$ws.Range("K113").Value = 1973.9333
$ws.Range("L113").Value = 74985.71000000001
$ws.Range("M113").Value = 196.0667000000001
$ws.Range("N113").Value = -79325.71000000001
# Row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 1672
$ws.Range("J113").Value = 4995
$ws.Range("L113").Value = 14985
$ws.Range("N113").Value = -19325
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 4791.8726
$ws.Range("I132").Value = 5496.5806
$ws.Range("J132").Value = 3426.5
$ws.Range("K132").Value = 16489.7418
$ws.Range("L132").Value = 10279.5
$ws.Range("M132").Value = -13959.7418
$ws.Range("N132").Value = -15339.5

